$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper data for the new handback entry: 2c9f4975-a5f1-4a10-8518-c03d0e449099
# ---------------------------------------------------------------------------
$newFile        = "2c9f4975-a5f1-4a10-8518-c03d0e449099.md"
$newFileDisp    = "e2e\2c9f4975-a5f1-4a10-8518-c03d0e449099.md"
$newStatus      = "Handed back: in sync with en-US"
$newOverviewDt  = "2016-09-07 08:42:16"

$oldFile        = "c500e649-6723-4740-b5a5-7740bf5d0c94.md"
$oldFileDisp    = "e2e\c500e649-6723-4740-b5a5-7740bf5d0c94.md"

# URLs (follow the existing naming / hashing convention used by the repo)
$newSrcUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c0983ea709d68e696bc102ece0c380dfa53c0e5/e2e/2c9f4975-a5f1-4a10-8518-c03d0e449099.md"
$newZhCnUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b8b66b7680d14445a539316678bafebaa97a8803/e2e/2c9f4975-a5f1-4a10-8518-c03d0e449099.md"
$newDeDeUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/df53e974a4d66bc25f9a3e8dc5daff54165b85c9/e2e/2c9f4975-a5f1-4a10-8518-c03d0e449099.md"

$oldSrcUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0ed3e565f6cc90b89e7b48debb05d97df36e3c7e/e2e/c500e649-6723-4740-b5a5-7740bf5d0c94.md"
$oldZhCnUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5b24d6d5480ddbb91b5c52a76a1f59e1282d5a73/e2e/c500e649-6723-4740-b5a5-7740bf5d0c94.md"
$oldDeDeUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/59df0355e2e4590f979408528ee31a7434ba771c/e2e/c500e649-6723-4740-b5a5-7740bf5d0c94.md"

# ===========================================================================
# Sheet 1: Overview
# ===========================================================================
$ws = $wb.Worksheets.Item("Overview")
$lo = $ws.ListObjects.Item(1)

# Insert a new data row before the existing "c500e649" row (row 3), which
# shifts it down to row 4 together with its values.
$ws.Rows.Item(3).Insert()
$lo.Resize($ws.Range("A1:G4"))

# Populate new row 3 (2c9f4975) -- mirrors the other "in sync" row.
$ws.Range("A3").Value = $newFile
$ws.Range("C3").Value = ".md"
$ws.Range("E3").Value = $newStatus
$ws.Range("F3").Value = $newStatus
$ws.Range("G3").Value = $newOverviewDt

# Fix up hyperlinks: remove the stale one anchored on row 3 (it still points
# at c500e649 even though that data moved to row 4), then recreate the two
# hyperlinks we need.
$hls = @($ws.Hyperlinks)
for ($i = $hls.Count - 1; $i -ge 0; $i--) {
    if ($hls[$i].Range.Row -eq 3) { $hls[$i].Delete() }
}
$ws.Hyperlinks.Add($ws.Range("B3"), $newSrcUrl, "", "", $newFileDisp)
$ws.Hyperlinks.Add($ws.Range("B4"), $oldSrcUrl, "", "", $oldFileDisp)

# ===========================================================================
# Sheet 2: zh-cn
# ===========================================================================
$ws = $wb.Worksheets.Item("zh-cn")
$lo = $ws.ListObjects.Item(1)

$ws.Rows.Item(3).Insert()
$lo.Resize($ws.Range("A1:P4"))

$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $newStatus
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = "2c9f4975-a5f1-4a10-8518-c03d0e449099.5560b1d4ff64036b1a2a928f24fb6387e121aaf0.zh-cn.xlf"
$ws.Range("H3").Value = "2016-09-07 08:41:56"
$ws.Range("J3").Value = "2c9f4975-a5f1-4a10-8518-c03d0e449099.5560b1d4ff64036b1a2a928f24fb6387e121aaf0.zh-cn.xlf"
$ws.Range("K3").Value = "2016-09-07 08:42:51"
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "True"
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "False"
$ws.Range("P3").Value = ""

$hls = @($ws.Hyperlinks)
for ($i = $hls.Count - 1; $i -ge 0; $i--) {
    if ($hls[$i].Range.Row -eq 3) { $hls[$i].Delete() }
}
$ws.Hyperlinks.Add($ws.Range("A3"), $newSrcUrl, "", "", $newFile)
$ws.Hyperlinks.Add($ws.Range("I3"), $newZhCnUrl, "", "", $newFile)
$ws.Hyperlinks.Add($ws.Range("A4"), $oldSrcUrl, "", "", $oldFile)
$ws.Hyperlinks.Add($ws.Range("I4"), $oldZhCnUrl, "", "", $oldFile)

# ===========================================================================
# Sheet 3: de-de
# ===========================================================================
$ws = $wb.Worksheets.Item("de-de")
$lo = $ws.ListObjects.Item(1)

$ws.Rows.Item(3).Insert()
$lo.Resize($ws.Range("A1:P4"))

$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $newStatus
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = "2c9f4975-a5f1-4a10-8518-c03d0e449099.5560b1d4ff64036b1a2a928f24fb6387e121aaf0.de-de.xlf"
$ws.Range("H3").Value = $newOverviewDt
$ws.Range("J3").Value = "2c9f4975-a5f1-4a10-8518-c03d0e449099.5560b1d4ff64036b1a2a928f24fb6387e121aaf0.de-de.xlf"
$ws.Range("K3").Value = "2016-09-07 08:43:16"
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "True"
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "False"
$ws.Range("P3").Value = ""

$hls = @($ws.Hyperlinks)
for ($i = $hls.Count - 1; $i -ge 0; $i--) {
    if ($hls[$i].Range.Row -eq 3) { $hls[$i].Delete() }
}
$ws.Hyperlinks.Add($ws.Range("A3"), $newSrcUrl, "", "", $newFile)
$ws.Hyperlinks.Add($ws.Range("I3"), $newDeDeUrl, "", "", $newFile)
$ws.Hyperlinks.Add($ws.Range("A4"), $oldSrcUrl, "", "", $oldFile)
$ws.Hyperlinks.Add($ws.Range("I4"), $oldDeDeUrl, "", "", $oldFile)
